$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.092.68'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.637.61'
$ws.Range('E3').Value = '  -1.96%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.10'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5261'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2600'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06315'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.73'
$ws.Range('E10').Value = '  -2.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07656'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.652.06'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.430'
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.15'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5507'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8192'
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.11'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.078.78'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.699'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.29'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.15'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.169'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.97'
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1218'
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.414'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.86'
$ws.Range('E28').Value = '  -0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.410'
$ws.Range('E29').Value = '  +4.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06015'
$ws.Range('E30').Value = '  -4.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.256'
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.413'
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.641'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9873'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.395'
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5743'
$ws.Range('E38').Value = '  -5.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01622'
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8549'
$ws.Range('E40').Value = '  -2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.040.60'
$ws.Range('E41').Value = '  -5.57%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.725'
$ws.Range('E43').Value = '  -6.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.65'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.788.14'
$ws.Range('E45').Value = '  -1.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -2.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.63'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.066'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4223'
$ws.Range('E51').Value = '  -0.56%  '
